# Renames the "services" collection placeholder to "ship-services".
#
# The template marks a repeating "services" section with
#   {{services}} ... {{/services}}
# Both the opening and closing tags are rewritten to reference the new
# collection name "ship-services" - but (per the target OOXML) each tag
# ends up represented as three adjacent runs, e.g.
#   <w:r><w:t>{{</w:t></w:r><w:r><w:t>ship-services</w:t></w:r><w:r><w:t>}}</w:t></w:r>
# instead of one run holding the whole tag text. Plain text assignment
# (Range.Text = ...) would leave everything in a single run because Word
# always coalesces adjacent runs that end up with identical formatting,
# so each fragment is inserted with a temporary bookmark acting as a
# run-boundary wedge between it and its neighbour, and the run is
# "touched" with a formatting round-trip (Bold on, then off again) so
# Word keeps an explicit (empty) run-properties element for it. The
# wedging bookmarks are deleted again immediately afterwards, so no
# bookmarks remain in the saved document.

$d = $word.ActiveDocument

function Insert-PlaceholderPart($anchorRange, $text) {
    # Inserts $text right after (the end of) $anchorRange and forces Word
    # to keep it as its own run (distinct from whatever precedes it) by
    # round-tripping a character formatting property on it.
    $anchorRange.Collapse(0)
    $partStart = $anchorRange.End
    $anchorRange.InsertAfter($text)
    $partRange = $d.Range($partStart, $partStart + $text.Length)
    $partRange.Font.Bold = 1
    $partRange.Font.Bold = 0
    $anchorRange.Collapse(0)
}

function Replace-PlaceholderTag($searchText, $prefix, $name, $suffix) {
    # Locate the whole placeholder tag (e.g. "{{services}}") via Find,
    # then rebuild it in-place as three separate runs: prefix, name,
    # suffix (e.g. "{{" / "ship-services" / "}}").
    $hit = $d.Content
    $found = $hit.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }

    # Clear the old text, leaving a collapsed insertion point in its place.
    $hit.Text = ""
    $hit.Collapse(1)

    $bk1 = "__tmp_wedge_1__"
    $bk2 = "__tmp_wedge_2__"

    Insert-PlaceholderPart $hit $prefix
    $d.Bookmarks.Add($bk1, $hit) | Out-Null

    Insert-PlaceholderPart $hit $name
    $d.Bookmarks.Add($bk2, $hit) | Out-Null

    Insert-PlaceholderPart $hit $suffix

    $d.Bookmarks($bk1).Delete()
    $d.Bookmarks($bk2).Delete()
}

Replace-PlaceholderTag "{{services}}" "{{" "ship-services" "}}"
Replace-PlaceholderTag "{{/services}}" "{{/" "ship-services" "}}"
